$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A, rows 2..97 currently hold "q1".."q96" and must become "q0".."q95"
# i.e. shift each quantile index down by one.
for ($r = 2; $r -le 97; $r++) {
    $n = $r - 2
    $ws.Cells.Item($r, 1).Value = "q$n"
}
